$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 479, shifting the existing rows 479:503 down to 480:504.
$ws.Rows("479:479").Insert()

# Populate the newly inserted row 479 with this week's data (same shape as
# the rows around it).
$ws.Range("A479").Value = 3
$ws.Range("B479").Value = "Femacal de La Calera"
$ws.Range("C479").Value = "Coquimbo"
$ws.Range("D479").Value = 45147
$ws.Range("E479").Value = 5
$ws.Range("F479").Value = 100112001
$ws.Range("G479").Value = "Berenjena"
$ws.Range("H479").Value = "Sin especificar"
$ws.Range("I479").Value = "Primera"
$ws.Range("J479").Value = 100
$ws.Range("K479").Value = 7000
$ws.Range("L479").Value = 7500
$ws.Range("M479").Value = 7275
$ws.Range("N479").Value = '$/caja 60 unidades'
$ws.Range("O479").Value = "Región de Arica y Parinacota"
$ws.Range("P479").Value = 121
$ws.Range("Q479").Value = 60
$ws.Range("R479").Value = "Hortaliza"
